# dynamicdata.xlsx - "updated 2-3 test cases"
#
# 1. Add two new sheets "InvalidLoginData" and "InvalidEmailData",
#    positioned right after "Sheet1" and right before "ContactUs".
# 2. Populate them with fake-credential / invalid-email test data.
# 3. Add a handful of extra (currently blank) formatted rows to the
#    bottom of "Sheet1" so the team has room to grow that table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the two new worksheets in the right order.
# Worksheets.Add() always inserts immediately before the ActiveSheet, so
# the freshly-added sheet lands right before "ContactUs" (the sheet that
# was active when the workbook was opened). We re-fetch each worksheet
# by name afterwards instead of reusing the variable returned by Add(),
# since that handle tracks a *position*, not the sheet object itself.
# ---------------------------------------------------------------------

$newSheet1 = $wb.Worksheets.Add()
$newSheet1.Name = "InvalidLoginData"

$newSheet2 = $wb.Worksheets.Add()
$newSheet2.Name = "InvalidEmailData"

# Fix tab order: InvalidLoginData, then InvalidEmailData, then ContactUs.
$invalidLoginData = $wb.Worksheets.Item("InvalidLoginData")
$invalidEmailData = $wb.Worksheets.Item("InvalidEmailData")
$invalidLoginData.Move($invalidEmailData)

# Re-fetch fresh references now that the sheets are in their final slots.
$invalidLoginData = $wb.Worksheets.Item("InvalidLoginData")
$invalidEmailData = $wb.Worksheets.Item("InvalidEmailData")

# ---------------------------------------------------------------------
# Step 2a: InvalidLoginData - bad credentials used by the login test.
# ---------------------------------------------------------------------

$invalidLoginData.Range("A1").Value = "email"
$invalidLoginData.Range("B1").Value = "password"
$invalidLoginData.Range("C1").Value = "expectedError"

$invalidLoginData.Range("A2").Value = "fakeuser@example.com"
$invalidLoginData.Range("B2").Value = "wrongPassword123"
$invalidLoginData.Range("C2").Value = "Your email or password is incorrect!"

$invalidLoginData.Hyperlinks.Add($invalidLoginData.Range("A2"), "mailto:fakeuser@example.com")

$invalidLoginData.Range("A1:C1").EntireColumn.AutoFit()

# ---------------------------------------------------------------------
# Step 2b: InvalidEmailData - malformed email used by the signup test.
# ---------------------------------------------------------------------

$invalidEmailData.Range("A1").Value = "name"
$invalidEmailData.Range("B1").Value = "invalidEmail"
$invalidEmailData.Range("C1").Value = "expectedValidation"

$invalidEmailData.Range("A2").Value = "John"
$invalidEmailData.Range("B2").Value = "plainaddress"
$invalidEmailData.Range("C2").Value = "include an '@'"

$invalidEmailData.Range("A2:C4").WrapText = $true
$invalidEmailData.Range("A2:C4").VerticalAlignment = -4108

$invalidEmailData.Activate()

# ---------------------------------------------------------------------
# Step 3: Sheet1 - extend the table with 5 additional blank rows,
# formatted like the rest of the data (rows 5-9, all 18 columns).
# ---------------------------------------------------------------------

$sheet1 = $wb.Worksheets.Item("Sheet1")
$blankRows = $sheet1.Range("A5:R9")
$blankRows.Style = "Normal 4"
$sheet1.Range("N5:R9").Select()
